$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.679.48"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.685.40"
$ws.Range("D3").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.678.81"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.33%  "

$ws.Range("E8").Value = "  +4.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.615"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.14"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.280.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "683.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.03"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.686.59"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +8.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.792.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.65%  "

$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.46%  "

$ws.Range("E21").Value = "  +3.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.945"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.94%  "

$ws.Range("E23").Value = "  +18.11%  "

$ws.Range("E24").Value = "  +4.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.21%  "

$ws.Range("E26").Value = "  +4.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.25%  "

$ws.Range("E28").Value = "  +5.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.46"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.90%  "

$ws.Range("E30").Value = "  +6.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +13.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.32"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "567.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.23%  "

$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.51"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.799.81"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.25%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("E39").Value = "  +5.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0467"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.22%  "

$ws.Range("E44").Value = "  +4.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.353"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.71%  "

$ws.Range("E46").Value = "  +9.04%  "

$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("E48").Value = "  +4.24%  "

$ws.Range("E49").Value = "  +3.21%  "

$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.90"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.99%  "
